$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (existing rows 2-16 shift down to 3-17)
$ws.Rows.Item(2).Insert()

# The inserted row picks up formatting from the header row (bold/bordered).
# Clear that so the new row matches the plain formatting used by the other data rows.
$ws.Range("A2:R2").ClearFormats()

# Fill in the new row 2 with the weekly record
$ws.Cells.Item(2, 1).Value = 10
$ws.Cells.Item(2, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(2, 3).Value = "La Araucanía"
$ws.Cells.Item(2, 4).Value = 44685
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value = 9
$ws.Cells.Item(2, 6).Value = 100112042
$ws.Cells.Item(2, 7).Value = "Locoto"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 60
$ws.Cells.Item(2, 11).Value = 5000
$ws.Cells.Item(2, 12).Value = 6000
$ws.Cells.Item(2, 13).Value = 5333
$ws.Cells.Item(2, 14).Value = "$/kilo"
$ws.Cells.Item(2, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(2, 16).Value = 5333
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = "Hortaliza"
